# edit.ps1
# Applies the "output generated at 456a3b4" update to 广州-漫展信息.xlsx:
#   - Refreshes the "want to go" head-count (column F) across all four sheets
#     (展览 / 演出 / 本地生活 / 全部类型) to match the latest scrape.
#   - Inserts a new event row into 全部类型 ("广州·浪漫古典·百年经典世界名曲音乐会",
#     2024-05-25) ahead of the existing "Look Look动漫嘉年华" row, shifting the
#     remaining rows down by one.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7661
$ws.Range("F3").Value = 97
$ws.Range("F4").Value = 77
$ws.Range("F5").Value = 6546
$ws.Range("F7").Value = 590
$ws.Range("F8").Value = 618
$ws.Range("F9").Value = 445
$ws.Range("F10").Value = 136
$ws.Range("F11").Value = 435
$ws.Range("F12").Value = 768
$ws.Range("F13").Value = 31
$ws.Range("F14").Value = 70
$ws.Range("F15").Value = 280
$ws.Range("F17").Value = 256
$ws.Range("F19").Value = 389
$ws.Range("F21").Value = 1084
$ws.Range("F23").Value = 589
$ws.Range("F24").Value = 2182
$ws.Range("F25").Value = 707
$ws.Range("F26").Value = 43
$ws.Range("F27").Value = 44
$ws.Range("F29").Value = 602
$ws.Range("F30").Value = 43

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 285
$ws.Range("F4").Value = 315
$ws.Range("F8").Value = 29

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 438

# --- Sheet 4: 全部类型 (value-only updates for rows 2-36) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 438
$ws.Range("F3").Value = 7661
$ws.Range("F4").Value = 97
$ws.Range("F5").Value = 77
$ws.Range("F6").Value = 285
$ws.Range("F7").Value = 6546
$ws.Range("F9").Value = 590
$ws.Range("F10").Value = 618
$ws.Range("F11").Value = 445
$ws.Range("F13").Value = 136
$ws.Range("F14").Value = 435
$ws.Range("F15").Value = 315
$ws.Range("F18").Value = 768
$ws.Range("F19").Value = 31
$ws.Range("F20").Value = 70
$ws.Range("F21").Value = 280
$ws.Range("F24").Value = 29
$ws.Range("F26").Value = 256
$ws.Range("F28").Value = 389
$ws.Range("F30").Value = 1084
$ws.Range("F32").Value = 589
$ws.Range("F33").Value = 2182
$ws.Range("F34").Value = 707
$ws.Range("F35").Value = 43
$ws.Range("F36").Value = 44

# --- Sheet 4: 全部类型 - insert new row for "广州·浪漫古典·百年经典世界名曲音乐会" ---
# Insert a blank row at position 39 (this pushes old row 39 "622排球少年only" down to row 40,
# and leaves old row 38 "Look Look动漫嘉年华" still at row 38 for now).
$ws.Rows.Item(39).Insert()

# Copy column-A style (bold/border/center) from row 38 into the newly inserted row 39's A cell,
# so it matches the formatting of all other index cells.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)

# Row 38 becomes the brand-new event (previously did not exist): 浪漫古典·百年经典世界名曲音乐会
# (Column B holds plain date-like text, not a real date value, in this workbook; prefix with an
# apostrophe to force text, then restore the default "Normal" style so no numeric/date formatting
# or quote-prefix marker sticks to the cell.)
$ws.Range("B38").Value = "'2024-05-25"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "广州·浪漫古典·百年经典世界名曲音乐会"
$ws.Range("D38").Value = "东风中路299号 广州中山纪念堂"
$ws.Range("E38").Value = "2024.05.25 20:00-05.25 21:30"
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 75
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=83327"
$ws.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202403/uRWx5ZEu1711079544682.jpeg"

# Row 39 (the newly inserted row) takes over what used to be row 38's event content
# (Look Look动漫嘉年华), with its updated want-to-go count.
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "'2024-06-01"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "广州·Look Look动漫嘉年华"
$ws.Range("D39").Value = "东沙大道16号 健康方舟"
$ws.Range("E39").Value = "2024.06.01 10:00-06.02 17:30"
$ws.Range("F39").Value = 602
$ws.Range("G39").Value = 29.9
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=82319"
$ws.Range("I39").Value = "//i2.hdslb.com/bfs/openplatform/202403/Zv7tuBjf1709620427087.png"

# Row 40 already holds the old row-39 event (622排球少年only) after the insert shifted it down;
# just refresh its index number and updated want-to-go count.
$ws.Range("A40").Value = 39
$ws.Range("F40").Value = 43

